$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# Trim the title down to its first run's text ("Lists") by deleting the
# remaining characters through a sub-range. This preserves the original
# (empty) run properties of the first run instead of fabricating a new
# run with synthesized properties.
$sub = $tr.Characters(6, $tr.Length - 5)
$sub.Text = ""

# Re-fetch the (now trimmed) text range and set the full desired text.
# Because "Lists" is a shared prefix, the existing first run is reused
# and simply extended, consolidating the paragraph into a single run.
$tr2 = $s.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "Lists (continued)"
